# --- Update the "Date" metadata value on the Metadata sheet (row 8, column B) ---
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-04-15T08:42:50-04:00"

# --- Add a new worksheet "Include ValueSets" placed right after "Include from LOINC" ---
$loincSheet = $wb.Worksheets.Item("Include from LOINC")
$newSheet = $wb.Worksheets.Add($null, $loincSheet)
$newSheet.Name = "Include ValueSets"

# Populate the header + value rows
$newSheet.Range("A1").Value = "ValueSet URL"
$newSheet.Range("A2").Value = "http://hl7.org/fhir/us/pacio-splasch/ValueSet/SPLASCHPureToneThresholdAudiometryPanel"

# Match the column widths used on the sibling "Include" sheet
$newSheet.Columns.Item(1).ColumnWidth = 29.8333333333333
$newSheet.Columns.Item(2).ColumnWidth = 49.8333333333333

# Re-use the existing bold-header / plain-value cell formats (copy formats only,
# so no new style entries get created and shared style indices are preserved)
$loincSheet.Range("A1").Copy() | Out-Null
$newSheet.Range("A1").PasteSpecial(-4122) | Out-Null
$loincSheet.Range("A2").Copy() | Out-Null
$newSheet.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the originally active sheet/tab selection
$meta.Activate()
